$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6438.2197
$ws.Range("J17").Value = 1599.175
$ws.Range("L17").Value = 4797.525
$ws.Range("N17").Value = -5133.525

$ws.Range("H38").Value = 965.6875
$ws.Range("I38").Value = 287.58334
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 862.7500200000001
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -490.7500200000001
$ws.Range("N38").Value = -9744

$ws.Range("H39").Value = 434.36365
$ws.Range("I39").Value = 30.583334
$ws.Range("J39").Value = 918.9
$ws.Range("K39").Value = 91.75000199999999
$ws.Range("L39").Value = 2756.7
$ws.Range("M39").Value = 204.249998
$ws.Range("N39").Value = -3348.7

$ws.Range("H86").Value = 13775.125
$ws.Range("I86").Value = 21440.6
$ws.Range("J86").Value = 999.3333
$ws.Range("K86").Value = 21440.6
$ws.Range("L86").Value = 999.3333
$ws.Range("M86").Value = -20317.6
$ws.Range("N86").Value = -3245.3333

$ws.Range("H89").Value = 13775.125
$ws.Range("I89").Value = 21440.6
$ws.Range("J89").Value = 999.3333
$ws.Range("K89").Value = 107203
$ws.Range("L89").Value = 4996.6665
$ws.Range("M89").Value = -101587
$ws.Range("N89").Value = -16228.6665

$ws.Range("H100").Value = 3231.5
$ws.Range("I100").Value = 2699.1875
$ws.Range("J100").Value = 7490
$ws.Range("K100").Value = 2699.1875
$ws.Range("L100").Value = 7490
$ws.Range("M100").Value = -2158.1875
$ws.Range("N100").Value = -8572

$ws.Range("H106").Value = 620.8
$ws.Range("I106").Value = 620.8
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 620.8
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = 10.20000000000005
$ws.Range("N106").ClearContents()

$ws.Range("H137").Value = 3192
$ws.Range("I137").Value = 3385.56
$ws.Range("K137").Value = 10156.68
$ws.Range("M137").Value = -7606.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 208301.5
$ws.Range("I6").Value = 356600
$ws.Range("J6").Value = 60003
$ws.Range("K6").Value = 356600
$ws.Range("L6").Value = 60003
$ws.Range("M6").Value = -356427
$ws.Range("N6").Value = -60349

$ws.Range("H26").Value = 1483.3334
$ws.Range("I26").Value = 975
$ws.Range("K26").Value = 975
$ws.Range("M26").Value = -645

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H39").Value = 4008
$ws.Range("I39").Value = 4008
$ws.Range("K39").Value = 4008
$ws.Range("M39").Value = -3488

$ws.Range("H45").Value = 1387.9412
$ws.Range("I45").Value = 1096.1538
$ws.Range("J45").Value = 2336.25
$ws.Range("K45").Value = 1096.1538
$ws.Range("L45").Value = 2336.25
$ws.Range("M45").Value = -719.1538
$ws.Range("N45").Value = -3090.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 15278.714
$ws.Range("I7").Value = 1251
$ws.Range("J7").Value = 25799.5
$ws.Range("K7").Value = 1251
$ws.Range("L7").Value = 25799.5
$ws.Range("M7").Value = -1138
$ws.Range("N7").Value = -26025.5

$ws.Range("H60").Value = 40000
$ws.Range("J60").Value = 40000
$ws.Range("L60").Value = 40000
$ws.Range("N60").Value = -41198

$ws.Range("H107").Value = 2726.2222
$ws.Range("I107").Value = 1748.3334
$ws.Range("J107").Value = 3704.111
$ws.Range("K107").Value = 1748.3334
$ws.Range("L107").Value = 3704.111
$ws.Range("M107").Value = 171.6666
$ws.Range("N107").Value = -7544.111

$ws.Range("H135").Value = 29700
$ws.Range("J135").Value = 29700
$ws.Range("L135").Value = 29700
$ws.Range("N135").Value = -39840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 17433.5
$ws.Range("I39").Value = 7413.3335
$ws.Range("J39").Value = 32463.75
$ws.Range("K39").Value = 7413.3335
$ws.Range("L39").Value = 32463.75
$ws.Range("M39").Value = -7022.3335
$ws.Range("N39").Value = -33245.75

$ws.Range("H49").Value = 17433.5
$ws.Range("I49").Value = 7413.3335
$ws.Range("J49").Value = 32463.75
$ws.Range("K49").Value = 7413.3335
$ws.Range("L49").Value = 32463.75
$ws.Range("M49").Value = -7231.3335
$ws.Range("N49").Value = -32827.75

$ws.Range("H52").Value = 29933.334
$ws.Range("J52").Value = 29933.334
$ws.Range("L52").Value = 29933.334
$ws.Range("N52").Value = -30521.334

$ws.Range("H127").Value = 32996.668
$ws.Range("J127").Value = 32996.668
$ws.Range("L127").Value = 32996.668
$ws.Range("N127").Value = -42916.668

$ws.Range("H141").Value = 29610
$ws.Range("J141").Value = 29610
$ws.Range("L141").Value = 29610
$ws.Range("N141").Value = -39970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 11666.777
$ws.Range("I6").Value = 166.83333
$ws.Range("J6").Value = 34666.668
$ws.Range("K6").Value = 500.49999
$ws.Range("L6").Value = 104000.004
$ws.Range("M6").Value = -387.49999
$ws.Range("N6").Value = -104226.004

$ws.Range("H36").Value = 1390.8
$ws.Range("I36").Value = 318
$ws.Range("J36").Value = 3000
$ws.Range("K36").Value = 954
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = -785
$ws.Range("N36").Value = -9338

$ws.Range("H50").Value = 101991.6
$ws.Range("I50").Value = 101
$ws.Range("J50").Value = 169918.67
$ws.Range("K50").Value = 303
$ws.Range("L50").Value = 509756.01
$ws.Range("M50").Value = 178
$ws.Range("N50").Value = -510718.01

$ws.Range("H53").Value = 101991.6
$ws.Range("I53").Value = 101
$ws.Range("J53").Value = 169918.67
$ws.Range("K53").Value = 303
$ws.Range("L53").Value = 509756.01
$ws.Range("M53").Value = 178
$ws.Range("N53").Value = -510718.01

$ws.Range("H59").Value = 1942.1428
$ws.Range("I59").Value = 231.66667
$ws.Range("J59").Value = 3225
$ws.Range("K59").Value = 695.00001
$ws.Range("L59").Value = 9675
$ws.Range("M59").Value = -155.00001
$ws.Range("N59").Value = -10755

$ws.Range("H64").Value = 2187.6
$ws.Range("J64").Value = 2755.5557
$ws.Range("L64").Value = 8266.667099999999
$ws.Range("N64").Value = -8806.667099999999

$ws.Range("H67").Value = 2187.6
$ws.Range("J67").Value = 2755.5557
$ws.Range("L67").Value = 8266.667099999999
$ws.Range("N67").Value = -10138.6671

$ws.Range("H75").Value = 2173.6
$ws.Range("J75").Value = 2947.5
$ws.Range("L75").Value = 8842.5
$ws.Range("N75").Value = -10838.5

$ws.Range("H78").Value = 2173.6
$ws.Range("J78").Value = 2947.5
$ws.Range("L78").Value = 26527.5
$ws.Range("N78").Value = -36511.5

$ws.Range("H87").Value = 11820
$ws.Range("J87").Value = 15816.667
$ws.Range("L87").Value = 47450.001
$ws.Range("N87").Value = -49946.001

$ws.Range("H90").Value = 11820
$ws.Range("J90").Value = 15816.667
$ws.Range("L90").Value = 142350.003
$ws.Range("N90").Value = -154830.003

$ws.Range("H131").Value = 1136.6383
$ws.Range("J131").Value = 1079.0286
$ws.Range("L131").Value = 3237.0858
$ws.Range("N131").Value = -13317.0858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12666.5
$ws.Range("I5").Value = 11000
$ws.Range("J5").Value = 12999.8
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12999.8
$ws.Range("M5").Value = -10888
$ws.Range("N5").Value = -13223.8

$ws.Range("H70").Value = 4736.5264
$ws.Range("I70").Value = 4654.923
$ws.Range("J70").Value = 4913.3335
$ws.Range("K70").Value = 4654.923
$ws.Range("L70").Value = 4913.3335
$ws.Range("M70").Value = -4384.923
$ws.Range("N70").Value = -5453.3335

$ws.Range("H73").Value = 4736.5264
$ws.Range("I73").Value = 4654.923
$ws.Range("J73").Value = 4913.3335
$ws.Range("K73").Value = 4654.923
$ws.Range("L73").Value = 4913.3335
$ws.Range("M73").Value = -3718.923
$ws.Range("N73").Value = -6785.3335

$ws.Range("H132").Value = 2935.7874
$ws.Range("I132").Value = 2429.818
$ws.Range("J132").Value = 4128.4287
$ws.Range("K132").Value = 7289.454000000001
$ws.Range("L132").Value = 12385.2861
$ws.Range("M132").Value = -4759.454000000001
$ws.Range("N132").Value = -17445.2861

$ws.Range("H137").Value = 29642.857
$ws.Range("J137").Value = 29642.857
$ws.Range("L137").Value = 29642.857
$ws.Range("N137").Value = -39842.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2266
$ws.Range("I96").Value = 2266
$ws.Range("K96").Value = 2266
$ws.Range("M96").Value = -893

$ws.Range("H100").Value = 988.4737
$ws.Range("I100").Value = 1017.875
$ws.Range("J100").Value = 967.0909
$ws.Range("K100").Value = 2035.75
$ws.Range("L100").Value = 1934.1818
$ws.Range("M100").Value = -1494.75
$ws.Range("N100").Value = -3016.1818

$ws.Range("H132").Value = 8564.914000000001
$ws.Range("I132").Value = 839.9216
$ws.Range("J132").Value = 29300.422
$ws.Range("K132").Value = 2519.7648
$ws.Range("L132").Value = 87901.266
$ws.Range("M132").Value = 10.23520000000008
$ws.Range("N132").Value = -92961.266
